$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 600, shifting existing rows 600:677 down to 601:678
$ws.Rows.Item(600).Insert()

# Populate the newly inserted row 600 with the new weekly record
$ws.Range("A600").Value = 8
$ws.Range("B600").Value = "Terminal La Palmera de La Serena"
$ws.Range("C600").Value = "Coquimbo"
$ws.Range("D600").Value = 45124
$ws.Range("E600").Value = 4
$ws.Range("F600").Value = 100112017
$ws.Range("G600").Value = "Apio"
$ws.Range("H600").Value = "Americana (o)"
$ws.Range("I600").Value = "Primera"
$ws.Range("J600").Value = 1200
$ws.Range("K600").Value = 7000
$ws.Range("L600").Value = 8000
$ws.Range("M600").Value = 7500
$ws.Range("N600").Value = "$/docena de matas"
$ws.Range("O600").Value = "Provincia del Elquí"
$ws.Range("P600").Value = 1250
$ws.Range("Q600").Value = 6
$ws.Range("R600").Value = "Hortaliza"
